$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 770.9231
$ws.Range("J41").Value = 1031.7142
$ws.Range("L41").Value = 1031.7142
$ws.Range("N41").Value = -1911.7142
$ws.Range("H43").Value = 877.7273
$ws.Range("I43").Value = 1233.3334
$ws.Range("J43").Value = 744.375
$ws.Range("K43").Value = 1233.3334
$ws.Range("L43").Value = 744.375
$ws.Range("M43").Value = -1164.3334
$ws.Range("N43").Value = -882.375
$ws.Range("H62").Value = 4569.4116
$ws.Range("I62").Value = 4636.25
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 4636.25
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -4012.25
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 4569.4116
$ws.Range("I65").Value = 4636.25
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 23181.25
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -20061.25
$ws.Range("N65").Value = -23740
$ws.Range("H107").Value = 328.6842
$ws.Range("I107").Value = 212.53847
$ws.Range("J107").Value = 580.3333
$ws.Range("K107").Value = 212.53847
$ws.Range("L107").Value = 580.3333
$ws.Range("M107").Value = 1707.46153
$ws.Range("N107").Value = -4420.3333
$ws.Range("H113").Value = 4156.2
$ws.Range("I113").Value = 3084.4443
$ws.Range("J113").Value = 4759.0625
$ws.Range("K113").Value = 3084.4443
$ws.Range("L113").Value = 4759.0625
$ws.Range("M113").Value = 169.5556999999999
$ws.Range("N113").Value = -11267.0625
$ws.Range("H132").Value = 9218.793
$ws.Range("I132").Value = 6746.5957
$ws.Range("J132").Value = 19781.818
$ws.Range("K132").Value = 20239.7871
$ws.Range("L132").Value = 59345.454
$ws.Range("M132").Value = -17709.7871
$ws.Range("N132").Value = -64405.454
$ws.Range("H137").Value = 15153732
$ws.Range("I137").Value = 26318282
$ws.Range("J137").Value = 1843.0714
$ws.Range("K137").Value = 78954846
$ws.Range("L137").Value = 5529.2142
$ws.Range("M137").Value = -78952296
$ws.Range("N137").Value = -10629.2142
$ws.Range("H138").Value = 1132.7894
$ws.Range("I138").Value = 913.1177
$ws.Range("K138").Value = 2739.3531
$ws.Range("M138").Value = 2400.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2589.1516
$ws.Range("I61").Value = 1575.5238
$ws.Range("J61").Value = 4363
$ws.Range("K61").Value = 1575.5238
$ws.Range("L61").Value = 4363
$ws.Range("M61").Value = -1363.5238
$ws.Range("N61").Value = -4787
$ws.Range("H136").Value = 2589.1516
$ws.Range("I136").Value = 1575.5238
$ws.Range("J136").Value = 4363
$ws.Range("K136").Value = 4726.5714
$ws.Range("L136").Value = 13089
$ws.Range("M136").Value = -2176.5714
$ws.Range("N136").Value = -18189

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1412.6072
$ws.Range("I107").Value = 1391.7916
$ws.Range("J107").Value = 1537.5
$ws.Range("K107").Value = 1391.7916
$ws.Range("L107").Value = 1537.5
$ws.Range("M107").Value = 528.2084
$ws.Range("N107").Value = -5377.5
$ws.Range("H134").Value = 4973.531
$ws.Range("I134").Value = 2442.6538
$ws.Range("J134").Value = 7834.522
$ws.Range("K134").Value = 7327.9614
$ws.Range("L134").Value = 23503.566
$ws.Range("M134").Value = -4792.9614
$ws.Range("N134").Value = -28573.566

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2177.125
$ws.Range("I31").Value = 1634.8823
$ws.Range("J31").Value = 5249.8335
$ws.Range("K31").Value = 1634.8823
$ws.Range("L31").Value = 5249.8335
$ws.Range("M31").Value = -1339.8823
$ws.Range("N31").Value = -5839.8335
$ws.Range("H34").Value = 2177.125
$ws.Range("I34").Value = 1634.8823
$ws.Range("J34").Value = 5249.8335
$ws.Range("K34").Value = 1634.8823
$ws.Range("L34").Value = 5249.8335
$ws.Range("M34").Value = -1432.8823
$ws.Range("N34").Value = -5653.8335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 8554.333000000001
$ws.Range("I23").Value = 96.333336
$ws.Range("J23").Value = 11373.667
$ws.Range("K23").Value = 289.000008
$ws.Range("L23").Value = 34121.001
$ws.Range("M23").Value = -54.00000799999998
$ws.Range("N23").Value = -34591.001
$ws.Range("H109").Value = 2667.1
$ws.Range("I109").Value = 2416.2
$ws.Range("J109").Value = 2918
$ws.Range("K109").Value = 7248.599999999999
$ws.Range("L109").Value = 8754
$ws.Range("M109").Value = -6208.599999999999
$ws.Range("N109").Value = -10834
$ws.Range("H113").Value = 563.53845
$ws.Range("I113").Value = 511.85715
$ws.Range("J113").Value = 623.8333
$ws.Range("K113").Value = 1535.57145
$ws.Range("L113").Value = 1871.4999
$ws.Range("M113").Value = 634.4285500000001
$ws.Range("N113").Value = -6211.4999
$ws.Range("H122").Value = 3498.5833
$ws.Range("I122").Value = 388.875
$ws.Range("J122").Value = 3781.2842
$ws.Range("K122").Value = 3499.875
$ws.Range("L122").Value = 34031.5578
$ws.Range("M122").Value = -1049.875
$ws.Range("N122").Value = -38931.5578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3759.8
$ws.Range("I113").Value = 3766.6667
$ws.Range("J113").Value = 3749.5
$ws.Range("K113").Value = 3766.6667
$ws.Range("L113").Value = 3749.5
$ws.Range("M113").Value = -1596.6667
$ws.Range("N113").Value = -8089.5
$ws.Range("H122").Value = 528867.75
$ws.Range("I122").Value = 910591.5600000001
$ws.Range("J122").Value = 3997.5
$ws.Range("K122").Value = 2731774.68
$ws.Range("L122").Value = 11992.5
$ws.Range("M122").Value = -2729324.68
$ws.Range("N122").Value = -16892.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 848.5625
$ws.Range("I22").Value = 833.2143
$ws.Range("J22").Value = 860.5
$ws.Range("K22").Value = 833.2143
$ws.Range("L22").Value = 860.5
$ws.Range("M22").Value = -538.2143
$ws.Range("N22").Value = -1450.5
$ws.Range("H27").Value = 848.5625
$ws.Range("I27").Value = 833.2143
$ws.Range("J27").Value = 860.5
$ws.Range("K27").Value = 833.2143
$ws.Range("L27").Value = 860.5
$ws.Range("M27").Value = -726.2143
$ws.Range("N27").Value = -1074.5
$ws.Range("H40").Value = 43481224
$ws.Range("J40").Value = 3584.0908
$ws.Range("L40").Value = 3584.0908
$ws.Range("N40").Value = -3856.0908
$ws.Range("H122").Value = 5295.5
$ws.Range("I122").Value = 6162.2085
$ws.Range("J122").Value = 4139.8887
$ws.Range("K122").Value = 18486.6255
$ws.Range("L122").Value = 12419.6661
$ws.Range("M122").Value = -16036.6255
$ws.Range("N122").Value = -17319.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5449.041
$ws.Range("I122").Value = 6100.9033
$ws.Range("J122").Value = 4326.3887
$ws.Range("K122").Value = 18302.7099
$ws.Range("L122").Value = 12979.1661
$ws.Range("M122").Value = -15852.7099
$ws.Range("N122").Value = -17879.1661
$ws.Range("H132").Value = 3324.4348
$ws.Range("I132").Value = 3346.4
$ws.Range("K132").Value = 10039.2
$ws.Range("M132").Value = -7509.200000000001
